$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (row 17): date 7/4/2023 (serial 45111), hours 0.5
$ws.Cells.Item(16, 1).Copy()
$ws.Cells.Item(17, 1).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Cells.Item(17, 1).Value = [DateTime]::FromOADate(45111)
$ws.Cells.Item(17, 2).Value = 0.5

# Update selection to reflect the new "next empty row" (A18)
$ws.Range("A18").Select()

$wb.Application.Calculate()
